# Auto-generated script to update crypto price/volume cells
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column (D) retains text storage, matching source data
# (values like '39.558.87' or '0.0₃0881' are not valid Excel numbers)
$ws.Range('D2:D51').NumberFormat = '@'

$ws.Range('D2').Value = '39.576.35'
$ws.Range('E2').Value = '  +1.63%  '
$ws.Range('D3').Value = '2.162.27'
$ws.Range('E3').Value = '  +0.39%  '
$ws.Range('D4').Value = '1.00'
$ws.Range('E4').Value = '  +0.10%  '
$ws.Range('D5').Value = '226.67'
$ws.Range('E5').Value = '  -0.83%  '
$ws.Range('D6').Value = '0.620'
$ws.Range('E6').Value = '  +0.12%  '
$ws.Range('D7').Value = '62.53'
$ws.Range('E7').Value = '  +0.28%  '
$ws.Range('E8').Value = '  +0.02%  '
$ws.Range('E9').Value = '  -0.49%  '
$ws.Range('D10').Value = '0.0846'
$ws.Range('E10').Value = '  -0.80%  '
$ws.Range('E11').Value = '  +0.57%  '
$ws.Range('D12').Value = '15.82'
$ws.Range('E12').Value = '  -1.17%  '
$ws.Range('D13').Value = '2.482.83'
$ws.Range('E13').Value = '  +0.61%  '
$ws.Range('D14').Value = '21.68'
$ws.Range('E14').Value = '  -2.41%  '
$ws.Range('D15').Value = '0.804'
$ws.Range('E15').Value = '  -1.56%  '
$ws.Range('E16').Value = '  -1.19%  '
$ws.Range('D17').Value = '2.157.73'
$ws.Range('E17').Value = '  +0.40%  '
$ws.Range('D18').Value = '39.568.40'
$ws.Range('E18').Value = '  +1.72%  '
$ws.Range('D19').Value = '71.56'
$ws.Range('E19').Value = '  -0.42%  '
$ws.Range('B20').Value = 'ShibaInu'
$ws.Range('C20').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D20').Value = '0.0₃0881'
$ws.Range('E20').Value = '  +3.68%  '
$ws.Range('B21').Value = 'Uniswap'
$ws.Range('C21').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D21').Value = '6.00'
$ws.Range('E21').Value = '  -2.17%  '
$ws.Range('D22').Value = '227.52'
$ws.Range('E22').Value = '  -0.06%  '
$ws.Range('E23').Value = '  +0.05%  '
$ws.Range('D24').Value = '2.34'
$ws.Range('E24').Value = '  +0.92%  '
$ws.Range('D25').Value = '2.31'
$ws.Range('E25').Value = '  -4.23%  '
$ws.Range('D26').Value = '170.30'
$ws.Range('E26').Value = '  -0.23%  '
$ws.Range('E27').Value = '  -3.07%  '
$ws.Range('E28').Value = '  +0.80%  '
$ws.Range('E29').Value = '  +1.86%  '
$ws.Range('D30').Value = '19.60'
$ws.Range('E30').Value = '  +0.12%  '
$ws.Range('D31').Value = '2.69'
$ws.Range('E31').Value = '  +4.57%  '
$ws.Range('E32').Value = '  +0.31%  '
$ws.Range('D33').Value = '4.46'
$ws.Range('E33').Value = '  -2.94%  '
$ws.Range('D34').Value = '4.69'
$ws.Range('E34').Value = '  -2.68%  '
$ws.Range('D35').Value = '6.95'
$ws.Range('E35').Value = '  -2.86%  '
$ws.Range('E36').Value = '  -0.04%  '
$ws.Range('D37').Value = '3.82'
$ws.Range('E37').Value = '  +7.79%  '
$ws.Range('E38').Value = '  -2.38%  '
$ws.Range('D39').Value = '1.00'
$ws.Range('E39').Value = '  +0.00%  '
$ws.Range('D40').Value = '4.92'
$ws.Range('E40').Value = '  +18.22%  '
$ws.Range('D41').Value = '102.11'
$ws.Range('E41').Value = '  -0.69%  '
$ws.Range('E42').Value = '  -1.60%  '
$ws.Range('D43').Value = '17.69'
$ws.Range('E43').Value = '  -2.62%  '
$ws.Range('D44').Value = '1.510.83'
$ws.Range('E44').Value = '  -1.56%  '
$ws.Range('E45').Value = '  +1.13%  '
$ws.Range('E46').Value = '  +0.54%  '
$ws.Range('E47').Value = '  +0.01%  '
$ws.Range('D48').Value = '0.0913'
$ws.Range('E48').Value = '  -0.55%  '
$ws.Range('E49').Value = '  -1.52%  '
$ws.Range('D50').Value = '0.000198'
$ws.Range('E50').Value = '  +35.07%  '
$ws.Range('E51').Value = '  -0.06%  '
